# "Added Backward extension option fo real-time data"
#
# The sheet holds a year-over-year forecast series ordered by date
# (column A, a date-formatted serial number), with columns B:E carrying
# paired (year, value) vectors. This adds 11 older observations
# (1984-1994) *before* the existing earliest row (1995), i.e. inserts
# 11 new rows right after the header row and shifts everything else
# down, then fills the freshly inserted rows with the backward-extended
# data. The worksheet dimension grows from A1:E31 to A1:E42.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 2
$insertCount = 11

# Insert 11 blank rows above the current row 2 (pushes 1995-2024 down to
# rows 13-42). Excel's row insert clones formatting from the row above
# (the bold/bordered header), so fix that up afterwards: columns B:E of
# data rows carry no explicit style, and column A carries the
# date-number-format style also used by every other data row (copy it
# from the row immediately below, which still has it after the shift).
$ws.Rows("$($firstNewRow):$($firstNewRow + $insertCount - 1)").Insert(-4121)
$ws.Range("B$($firstNewRow):E$($firstNewRow + $insertCount - 1)").ClearFormats()
$ws.Range("A$($firstNewRow + $insertCount)").Copy()
$ws.Range("A$($firstNewRow):A$($firstNewRow + $insertCount - 1)").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# date-serial, y_0 year, y_0 value, y_1 year, y_1 value
$newData = @(
    @(31047, 1984, 2.681929770019686,  1985, 2.73475129348093),
    @(31412, 1985, 2.573947575822677,  1986, 2.58358492564803),
    @(31777, 1986, 2.161034240664228,  1987, 2.388724974429235),
    @(32142, 1987, 0.9697275934645422, 1988, 2.698368179641242),
    @(32508, 1988, 3.194969449935003,  1989, 2.654510774528207),
    @(32873, 1989, 3.85009945173751,   1990, 2.342799083309055),
    @(33238, 1990, 5.073362306219398,  1991, 2.977303796668029),
    @(33603, 1991, 6.091605135014255,  1992, 2.066726874661873),
    @(33969, 1992, 2.064701871240571,  1993, 2.112386427028046),
    @(34334, 1993, -1.000531514043412, 1994, 2.575999544954621),
    @(34699, 1994, 2.998503002360153,  1995, 2.954478109176528)
)

for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $firstNewRow + $i
    $vals = $newData[$i]
    $ws.Cells.Item($r, 1).Value2 = $vals[0]
    $ws.Cells.Item($r, 2).Value2 = $vals[1]
    $ws.Cells.Item($r, 3).Value2 = $vals[2]
    $ws.Cells.Item($r, 4).Value2 = $vals[3]
    $ws.Cells.Item($r, 5).Value2 = $vals[4]
}

"Inserted $insertCount backward-extension rows; dimension now A1:E$($ws.Range('A1').End(-4121).Row)"
